# Auto-generated edit script applying scheduled-runner profit-sheet updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2228497.2
$ws.Range("J17").Value = 2259882.8
$ws.Range("L17").Value = 6779648.399999999
$ws.Range("N17").Value = -6779984.399999999

$ws.Range("H45").Value = 23750
$ws.Range("J45").Value = 23750
$ws.Range("L45").Value = 71250
$ws.Range("N45").Value = -71634

$ws.Range("H55").Value = 385.53845
$ws.Range("I55").Value = 999
$ws.Range("J55").Value = 201.5
$ws.Range("K55").Value = 999
$ws.Range("L55").Value = 201.5
$ws.Range("M55").Value = -785
$ws.Range("N55").Value = -629.5

$ws.Range("H70").Value = 2308.3667
$ws.Range("I70").Value = 3800.5
$ws.Range("J70").Value = 2078.8076
$ws.Range("K70").Value = 11401.5
$ws.Range("L70").Value = 6236.4228
$ws.Range("M70").Value = -11131.5
$ws.Range("N70").Value = -6776.4228

$ws.Range("H73").Value = 2308.3667
$ws.Range("I73").Value = 3800.5
$ws.Range("J73").Value = 2078.8076
$ws.Range("K73").Value = 11401.5
$ws.Range("L73").Value = 6236.4228
$ws.Range("M73").Value = -10465.5
$ws.Range("N73").Value = -8108.4228

$ws.Range("H106").Value = 60608560
$ws.Range("I106").Value = 23812024
$ws.Range("J106").Value = 125002500
$ws.Range("K106").Value = 23812024
$ws.Range("L106").Value = 125002500
$ws.Range("M106").Value = -23811393
$ws.Range("N106").Value = -125003762

$ws.Range("H132").Value = 7752966.5
$ws.Range("I132").Value = 782.6177
$ws.Range("J132").Value = 37038996
$ws.Range("K132").Value = 2347.8531
$ws.Range("L132").Value = 111116988
$ws.Range("M132").Value = 182.1468999999997
$ws.Range("N132").Value = -111122048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1362.8572
$ws.Range("I2").Value = 250.5
$ws.Range("J2").Value = 1807.8
$ws.Range("K2").Value = 250.5
$ws.Range("L2").Value = 1807.8
$ws.Range("M2").Value = -137.5
$ws.Range("N2").Value = -2033.8

$ws.Range("H32").Value = 6661.3857
$ws.Range("I32").Value = 4971.6665
$ws.Range("K32").Value = 4971.6665
$ws.Range("M32").Value = -4684.6665

$ws.Range("H61").Value = 4046.775
$ws.Range("I61").Value = 4511.5938
$ws.Range("K61").Value = 4511.5938
$ws.Range("M61").Value = -4299.5938

$ws.Range("H102").Value = 1951368
$ws.Range("I102").Value = 2180647
$ws.Range("K102").Value = 2180647
$ws.Range("M102").Value = -2179025

$ws.Range("H116").Value = 1362.8572
$ws.Range("I116").Value = 250.5
$ws.Range("J116").Value = 1807.8
$ws.Range("K116").Value = 250.5
$ws.Range("L116").Value = 1807.8
$ws.Range("M116").Value = 2043.5
$ws.Range("N116").Value = -6395.8

$ws.Range("H136").Value = 4046.775
$ws.Range("I136").Value = 4511.5938
$ws.Range("K136").Value = 13534.7814
$ws.Range("M136").Value = -10984.7814

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1362.8572
$ws.Range("I3").Value = 250.5
$ws.Range("J3").Value = 1807.8
$ws.Range("K3").Value = 250.5
$ws.Range("L3").Value = 1807.8
$ws.Range("M3").Value = -136.5
$ws.Range("N3").Value = -2035.8

$ws.Range("H134").Value = 4593.3335
$ws.Range("I134").Value = 5775.926
$ws.Range("J134").Value = 2464.6667
$ws.Range("K134").Value = 17327.778
$ws.Range("L134").Value = 7394.000100000001
$ws.Range("M134").Value = -14792.778
$ws.Range("N134").Value = -12464.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5996466
$ws.Range("I31").Value = 1571.26
$ws.Range("J31").Value = 14322708
$ws.Range("K31").Value = 1571.26
$ws.Range("L31").Value = 14322708
$ws.Range("M31").Value = -1276.26
$ws.Range("N31").Value = -14323298

$ws.Range("H34").Value = 5996466
$ws.Range("I34").Value = 1571.26
$ws.Range("J34").Value = 14322708
$ws.Range("K34").Value = 1571.26
$ws.Range("L34").Value = 14322708
$ws.Range("M34").Value = -1369.26
$ws.Range("N34").Value = -14323112

$ws.Range("H58").Value = 2646449
$ws.Range("I58").Value = 3333866.8
$ws.Range("K58").Value = 3333866.8
$ws.Range("M58").Value = -3333663.8

$ws.Range("H70").Value = 40000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null

$ws.Range("H73").Value = 40000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null

$ws.Range("H107").Value = 12346236
$ws.Range("I107").Value = 16667030
$ws.Range("J107").Value = 1110.1428
$ws.Range("K107").Value = 16667030
$ws.Range("L107").Value = 1110.1428
$ws.Range("M107").Value = -16665110
$ws.Range("N107").Value = -4950.1428

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null

$ws.Range("H132").Value = 2565377
$ws.Range("I132").Value = 3334424
$ws.Range("J132").Value = 1887.6111
$ws.Range("K132").Value = 10003272
$ws.Range("L132").Value = 5662.8333
$ws.Range("M132").Value = -10000742
$ws.Range("N132").Value = -10722.8333

$ws.Range("H134").Value = 4976655
$ws.Range("I134").Value = 8548896
$ws.Range("J134").Value = 1034.3214
$ws.Range("K134").Value = 25646688
$ws.Range("L134").Value = 3102.9642
$ws.Range("M134").Value = -25644153
$ws.Range("N134").Value = -8172.9642

$ws.Range("H136").Value = 2646449
$ws.Range("I136").Value = 3333866.8
$ws.Range("K136").Value = 10001600.4
$ws.Range("M136").Value = -9999050.399999999

$ws.Range("H141").Value = 33294.777
$ws.Range("J141").Value = 33294.777
$ws.Range("L141").Value = 33294.777
$ws.Range("N141").Value = -43654.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 28572650
$ws.Range("I75").Value = 313
$ws.Range("J75").Value = 35715732
$ws.Range("K75").Value = 939
$ws.Range("L75").Value = 107147196
$ws.Range("M75").Value = 59
$ws.Range("N75").Value = -107149192

$ws.Range("H76").Value = 3625
$ws.Range("J76").Value = 3625
$ws.Range("L76").Value = 10875
$ws.Range("N76").Value = -11641

$ws.Range("H78").Value = 28572650
$ws.Range("I78").Value = 313
$ws.Range("J78").Value = 35715732
$ws.Range("K78").Value = 2817
$ws.Range("L78").Value = 321441588
$ws.Range("M78").Value = 2175
$ws.Range("N78").Value = -321451572

$ws.Range("H79").Value = 3625
$ws.Range("J79").Value = 3625
$ws.Range("L79").Value = 10875
$ws.Range("N79").Value = -13527

$ws.Range("H113").Value = 2609265
$ws.Range("J113").Value = 909616.9399999999
$ws.Range("L113").Value = 2728850.82
$ws.Range("N113").Value = -2733190.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5644.95
$ws.Range("I70").Value = 5278.5713
$ws.Range("K70").Value = 5278.5713
$ws.Range("M70").Value = -5008.5713

$ws.Range("H73").Value = 5644.95
$ws.Range("I73").Value = 5278.5713
$ws.Range("K73").Value = 5278.5713
$ws.Range("M73").Value = -4342.5713

$ws.Range("H102").Value = 1606.5834
$ws.Range("I102").Value = 1327.8334
$ws.Range("J102").Value = 1885.3334
$ws.Range("K102").Value = 1327.8334
$ws.Range("L102").Value = 1885.3334
$ws.Range("M102").Value = 294.1666
$ws.Range("N102").Value = -5129.3334

$ws.Range("H113").Value = 71430110
$ws.Range("I113").Value = 142858340
$ws.Range("J113").Value = 1885.7142
$ws.Range("K113").Value = 142858340
$ws.Range("L113").Value = 1885.7142
$ws.Range("M113").Value = -142856170
$ws.Range("N113").Value = -6225.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null

$ws.Range("H132").Value = 13012431
$ws.Range("I132").Value = 14819185
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 44457555
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -44455025
$ws.Range("N132").Value = -16460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3116.6667
$ws.Range("I62").Value = 3060
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 3060
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -2436
$ws.Range("N62").Value = -4648

$ws.Range("H65").Value = 3116.6667
$ws.Range("I65").Value = 3060
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 15300
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -12180
$ws.Range("N65").Value = -23240

$ws.Range("H113").Value = 1027.1143
$ws.Range("I113").Value = 994.4483
$ws.Range("J113").Value = 1185
$ws.Range("K113").Value = 2983.3449
$ws.Range("L113").Value = 3555
$ws.Range("M113").Value = -813.3449000000001
$ws.Range("N113").Value = -7895

$ws.Range("H126").Value = 1288.7858
$ws.Range("I126").Value = 895.6
$ws.Range("J126").Value = 1507.2222
$ws.Range("K126").Value = 2686.8
$ws.Range("L126").Value = 4521.6666
$ws.Range("M126").Value = -216.8000000000002
$ws.Range("N126").Value = -9461.6666
